# Status update as of 4-may-2017 (PILOT_Status workbook, "files" sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("files")

# ds.xpt and ex.xpt are now done (were InProgress); ex.xpt's reviewer guide note is now "yes"
$ws.Range("C5").Value = "done"
$ws.Range("C6").Value = "done"
$ws.Range("E6").Value = "yes"

# lb.xpt and supplb.xpt are now assigned InProgress to Cindy
$ws.Range("C7").Value = "InProgress"
$ws.Range("D7").Value = "Cindy"
$ws.Range("C16").Value = "InProgress"
$ws.Range("D16").Value = "Cindy"

# notes column: qs.xpt gets a new note, relrec.xpt's note is corrected/expanded
$ws.Range("F9").Value = "fix dataset name in xpt file"
$ws.Range("F10").Value = "has references to DS and AE"

# suppds.xpt is now done
$ws.Range("C15").Value = "done"

# new "Assigned To" header column
$ws.Range("D1").Value = "Assigned To"

$ws.Range("E11").Select()
